# Refresh the cryptos list: updated prices / 1h volume % for rows 2-50,
# and row 51 swapped from BabyDogeCoin to Cronos (name, link, price, volume%).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.793.75"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.605.67"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'212.98"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "'28.23"
$ws.Range("E8").Value = "  +5.15%  "
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("D11").Value = "'0.0909"
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").Value = "1.836.80"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Value = "1.613.02"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("E14").Value = "  +3.78%  "
$ws.Range("D15").Value = "29.760.22"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").Value = "'64.07"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "'242.40"
$ws.Range("E18").Value = "  -2.27%  "
$ws.Range("D19").Value = "'7.89"
$ws.Range("E19").Value = "  +3.43%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("D23").Value = "'9.39"
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("D25").Value = "'155.22"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("D26").Value = "'15.47"
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("D32").Value = "'3.24"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("E33").Value = "  +2.53%  "
$ws.Range("D34").Value = "1.420.71"
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("E35").Value = "  +2.84%  "
$ws.Range("E36").Value = "  +2.36%  "
$ws.Range("E37").Value = "  -2.01%  "
$ws.Range("D38").Value = "'2.29"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D40").Value = "'0.547"
$ws.Range("E40").Value = "  +1.89%  "
$ws.Range("D41").Value = "'57.15"
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("D42").Value = "'0.0496"
$ws.Range("E42").Value = "  +5.93%  "
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("E44").Value = "  +1.83%  "
$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("D46").Value = "'67.01"
$ws.Range("E46").Value = "  -2.34%  "
$ws.Range("E47").Value = "  +17.51%  "
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("D49").Value = "1.745.42"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").Value = "'86.51"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0524"
$ws.Range("E51").Value = "  +0.68%  "
